$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column B (shifts old B/C/D -> C/D/E), making room for "_requirements"
$ws.Columns.Item(2).Insert()

# 2. Insert a new row before row 4 (the old "enasarco" row, now at row 5 after the column
#    insert didn't change row count) so we can add the G=='zero' split for wt_1038.
$ws.Rows.Item(4).Insert()

# 3. Append two more blank rows at the bottom for the new enasarco "not zero" row and the
#    brand-new z0bug.wt_1040-23A_1 row.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(7).Insert()

# --- Header row ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "_requirements"
$ws.Range("C1").Value = "withholding_tax_id"
$ws.Range("D1").Value = "tax"
$ws.Range("E1").Value = "base"

# --- Row 2: z0bug.wt_1040_1 / z0bug.wt_1040 (unchanged values, just shifted right) ---
$ws.Range("A2").Value = "z0bug.wt_1040_1"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "z0bug.wt_1040"
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 1

# --- Row 3: z0bug.wt_1038_1 / z0bug.wt_1038, requirement G=='zero' ---
$ws.Range("A3").Value = "z0bug.wt_1038_1"
$ws.Range("B3").Value = "G=='zero'"
$ws.Range("C3").Value = "z0bug.wt_1038"
$ws.Range("D3").Value = 23
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.5"
$ws.Range("E3").NumberFormat = "General"

# --- Row 4 (new): z0bug.wt_1038_1 / z0bug.wt_1038, requirement G!='zero' ---
$ws.Range("A4").Value = "z0bug.wt_1038_1"
$ws.Range("B4").Value = "G!='zero'"
$ws.Range("C4").Value = "z0bug.wt_1038"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "11.5"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = 1

# --- Row 5: z0bug.wt_enasarco_1_1 / z0bug.wt_enasarco_1, requirement G=='zero' ---
$ws.Range("A5").Value = "z0bug.wt_enasarco_1_1"
$ws.Range("B5").Value = "G=='zero'"
$ws.Range("C5").Value = "z0bug.wt_enasarco_1"
$ws.Range("D5").Value = 17
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.5"
$ws.Range("E5").NumberFormat = "General"

# --- Row 6 (new): z0bug.wt_enasarco_1_1 / z0bug.wt_enasarco_1, requirement G!='zero' ---
$ws.Range("A6").Value = "z0bug.wt_enasarco_1_1"
$ws.Range("B6").Value = "G!='zero'"
$ws.Range("C6").Value = "z0bug.wt_enasarco_1"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.5"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = 1

# --- Row 7 (new): z0bug.wt_1040-23A_1 / z0bug.wt_1040-23A (no requirement) ---
$ws.Range("A7").Value = "z0bug.wt_1040-23A_1"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "z0bug.wt_1040-23A"
$ws.Range("D7").Value = 23
$ws.Range("E7").Value = 1

# --- Column widths ---
# (The stored OOXML "width" is ColumnWidth + ~0.8333 char-units of built-in
# padding, then the host quantizes to the nearest 1/6 of a character; these
# inputs are calibrated to land as close as possible to the authored widths
# of 21.02 / 14.88 / 19.58 / 4.07 / 5.46.)
$ws.Columns.Item(1).ColumnWidth = 20.186666666666667
$ws.Columns.Item(2).ColumnWidth = 14.046666666666667
$ws.Columns.Item(3).ColumnWidth = 18.746666666666666
$ws.Columns.Item(4).ColumnWidth = 3.236666666666667
$ws.Columns.Item(5).ColumnWidth = 4.626666666666667

# --- Selection, matching the authored selection state ---
$ws.Range("D5").Select()
